# Edit fig on pred result and caption of it
# Adjusts the timing-label textboxes (re-aligning + relabeling) and
# shifts the (a)/(b)/(c) caption textboxes on slide 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$EMU_PER_PT = 12700

function Set-ShapePos {
    param($shape, $xEmu, $yEmu)
    # +0.5 EMU nudges the EMU->point conversion so the host's
    # single-precision round-trip lands back on the exact target EMU
    # (the host truncates rather than rounds on point->EMU write-back).
    $shape.Left = ($xEmu + 0.5) / 12700
    $shape.Top  = ($yEmu + 0.5) / 12700
}

# "1 second" label (position only, text unchanged)
$tb27 = $s.Shapes.Item("TextBox 27")
Set-ShapePos $tb27 2600307 921253

# "3 seconds" -> "2 seconds"
$tb17 = $s.Shapes.Item("TextBox 17")
Set-ShapePos $tb17 6414101 917455
$tb17.TextFrame.TextRange.Text = "2 seconds"

# "5 seconds" -> "3 seconds"
$tb18 = $s.Shapes.Item("TextBox 18")
Set-ShapePos $tb18 10245235 921253
$tb18.TextFrame.TextRange.Text = "3 seconds"

# "7 seconds" -> "4 seconds"
$tb19 = $s.Shapes.Item("TextBox 19")
Set-ShapePos $tb19 14112571 912425
$tb19.TextFrame.TextRange.Text = "4 seconds"

# Caption "(a)" - shift horizontally only
$tb28 = $s.Shapes.Item("TextBox 28")
Set-ShapePos $tb28 814699 2683846

# Caption "(b)" - shift horizontally only
$tb29 = $s.Shapes.Item("TextBox 29")
Set-ShapePos $tb29 814699 5878457

# Caption "(c)" - shift horizontally only
$tb30 = $s.Shapes.Item("TextBox 30")
Set-ShapePos $tb30 814699 8425157
